# Scheduled-runner style refresh of the FFXIV Leve-profit market data.
# Updates the currentAveragePrice*/LevePrice*/LeveProfit* columns (H, I, J, K,
# L, M, N) on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with
# freshly pulled marketboard figures. Cells that no longer have a computed
# value (e.g. because a dependent price dropped out) are cleared instead of
# being left with a stale number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 248.13637
$ws.Range("I33").Value = 257.57144
$ws.Range("K33").Value = 257.57144
$ws.Range("M33").Value = -28.57144
$ws.Range("H38").Value = 1335.1428
$ws.Range("I38").Value = 472.83334
$ws.Range("K38").Value = 1418.50002
$ws.Range("M38").Value = -1046.50002
$ws.Range("H58").Value = 8369.5
$ws.Range("J58").Value = 12354.25
$ws.Range("L58").Value = 37062.75
$ws.Range("N58").Value = -37362.75
$ws.Range("H70").Value = 1429454.1
$ws.Range("J70").Value = 2500745
$ws.Range("L70").Value = 7502235
$ws.Range("N70").Value = -7502775
$ws.Range("H73").Value = 1429454.1
$ws.Range("J73").Value = 2500745
$ws.Range("L73").Value = 7502235
$ws.Range("N73").Value = -7504107
$ws.Range("H80").Value = 661.75
$ws.Range("I80").Value = 677.2
$ws.Range("J80").Value = 654.7273
$ws.Range("K80").Value = 2031.6
$ws.Range("L80").Value = 1964.1819
$ws.Range("M80").Value = -1033.6
$ws.Range("N80").Value = -3960.1819
$ws.Range("H83").Value = 661.75
$ws.Range("I83").Value = 677.2
$ws.Range("J83").Value = 654.7273
$ws.Range("K83").Value = 6094.8
$ws.Range("L83").Value = 5892.545700000001
$ws.Range("M83").Value = -1102.8
$ws.Range("N83").Value = -15876.5457
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4757.3696
$ws.Range("J32").Value = 14666.333
$ws.Range("L32").Value = 14666.333
$ws.Range("N32").Value = -15240.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 824.5
$ws.Range("I22").Value = 824.5
$ws.Range("K22").Value = 824.5
$ws.Range("M22").Value = -651.5
$ws.Range("H64").Value = 1389.6
$ws.Range("J64").Value = 1487
$ws.Range("L64").Value = 1487
$ws.Range("N64").Value = -1937
$ws.Range("H67").Value = 1389.6
$ws.Range("J67").Value = 1487
$ws.Range("L67").Value = 1487
$ws.Range("N67").Value = -3047
$ws.Range("H99").Value = 930.25
$ws.Range("J99").Value = 910.5
$ws.Range("L99").Value = 910.5
$ws.Range("N99").Value = -3906.5
$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 59000
$ws.Range("J53").Value = 59000
$ws.Range("L53").Value = 59000
$ws.Range("N53").Value = -60214
$ws.Range("H132").Value = 4269.1665
$ws.Range("I132").Value = 4211.1665
$ws.Range("J132").Value = 4385.1665
$ws.Range("K132").Value = 12633.4995
$ws.Range("L132").Value = 13155.4995
$ws.Range("M132").Value = -10103.4995
$ws.Range("N132").Value = -18215.4995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 81.25
$ws.Range("J2").Value = 87
$ws.Range("L2").Value = 522
$ws.Range("N2").Value = -748
$ws.Range("H6").Value = 974
$ws.Range("I6").Value = 298.66666
$ws.Range("K6").Value = 895.9999799999999
$ws.Range("M6").Value = -782.9999799999999
$ws.Range("H7").Value = 232.75
$ws.Range("I7").Value = 232.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 698.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -586.25
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 1500
$ws.Range("M8").Value = -1361
$ws.Range("H17").Value = 1100
$ws.Range("J17").Value = 1352.5
$ws.Range("L17").Value = 4057.5
$ws.Range("N17").Value = -4395.5
$ws.Range("H23").Value = 214.11111
$ws.Range("J23").Value = 214.11111
$ws.Range("L23").Value = 642.3333299999999
$ws.Range("N23").Value = -1112.33333
$ws.Range("H34").Value = 2774.5
$ws.Range("J34").Value = 7000
$ws.Range("L34").Value = 21000
$ws.Range("N34").Value = -21168
$ws.Range("H55").Value = 6838.8
$ws.Range("J55").Value = 15000
$ws.Range("L55").Value = 45000
$ws.Range("N55").Value = -45354
$ws.Range("H113").Value = 2642
$ws.Range("J113").Value = 2642
$ws.Range("L113").Value = 7926
$ws.Range("N113").Value = -12266
$ws.Range("H122").Value = 92544.45
$ws.Range("J122").Value = 144571
$ws.Range("L122").Value = 1301139
$ws.Range("N122").Value = -1306039
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -378
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 574.6
$ws.Range("I107").Value = 518.25
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 518.25
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1401.75
$ws.Range("N107").Value = -4640
$ws.Range("H113").Value = 1950
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 220
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 1479.4
$ws.Range("I126").Value = 1483.5
$ws.Range("J126").Value = 1476.6666
$ws.Range("K126").Value = 4450.5
$ws.Range("L126").Value = 4429.9998
$ws.Range("M126").Value = -1980.5
$ws.Range("N126").Value = -9369.9998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3555.1
$ws.Range("I61").Value = 3583.3333
$ws.Range("J61").Value = 3512.75
$ws.Range("K61").Value = 3583.3333
$ws.Range("L61").Value = 3512.75
$ws.Range("M61").Value = -3381.3333
$ws.Range("N61").Value = -3916.75
$ws.Range("H82").Value = 1844.5
$ws.Range("I82").Value = 1721
$ws.Range("J82").Value = 2462
$ws.Range("K82").Value = 1721
$ws.Range("L82").Value = 2462
$ws.Range("M82").Value = -1360
$ws.Range("N82").Value = -3184
$ws.Range("H85").Value = 1844.5
$ws.Range("I85").Value = 1721
$ws.Range("J85").Value = 2462
$ws.Range("K85").Value = 1721
$ws.Range("L85").Value = 2462
$ws.Range("M85").Value = -473
$ws.Range("N85").Value = -4958
$ws.Range("H100").Value = 1331.6666
$ws.Range("I100").Value = 497.5
$ws.Range("K100").Value = 497.5
$ws.Range("M100").Value = 43.5
$ws.Range("H113").Value = 3555.1
$ws.Range("I113").Value = 3583.3333
$ws.Range("J113").Value = 3512.75
$ws.Range("K113").Value = 3583.3333
$ws.Range("L113").Value = 3512.75
$ws.Range("M113").Value = -1413.3333
$ws.Range("N113").Value = -7852.75
$ws.Range("H122").Value = 8166
$ws.Range("I122").Value = 7249.5
$ws.Range("K122").Value = 21748.5
$ws.Range("M122").Value = -19298.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4643.5625
$ws.Range("I81").Value = 2784.8462
$ws.Range("K81").Value = 5569.6924
$ws.Range("M81").Value = -4508.6924
$ws.Range("H84").Value = 4643.5625
$ws.Range("I84").Value = 2784.8462
$ws.Range("K84").Value = 27848.462
$ws.Range("M84").Value = -22544.462
$ws.Range("H122").Value = 4566.3335
$ws.Range("I122").Value = 4479.6
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 13438.8
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10988.8
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 2586.5557
$ws.Range("I126").Value = 2534.875
$ws.Range("K126").Value = 7604.625
$ws.Range("M126").Value = -5134.625
